$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (row -> D, L, M, N, O, P, S) after the shuffle/update.
# Only columns D, L, M, N, O, P, S change; everything else stays the same.
$data = @{
    2  = @(44511, "Primera", 120, 28000, 28000, 28000, 2800)
    3  = @(44434, "Primera", 20,  20000, 20000, 20000, 2000)
    4  = @(44517, "Especial", 100, 27000, 27000, 27000, 2700)
    5  = @(44517, "Primera", 30,  25000, 25000, 25000, 2500)
    6  = @(44473, "Primera", 180, 20000, 20000, 20000, 2000)
    7  = @(44435, "Primera", 40,  20000, 20000, 20000, 2000)
    8  = @(44476, "Primera", 120, 20000, 20000, 20000, 2000)
    9  = @(44503, "Primera", 60,  30000, 30000, 30000, 3000)
    10 = @(44503, "Segunda", 50,  25000, 25000, 25000, 2500)
    11 = @(44466, "Primera", 60,  20000, 20000, 20000, 2000)
    12 = @(44432, "Primera", 20,  20000, 20000, 20000, 2000)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("L$row").Value = $vals[1]
    $ws.Range("M$row").Value = $vals[2]
    $ws.Range("N$row").Value = $vals[3]
    $ws.Range("O$row").Value = $vals[4]
    $ws.Range("P$row").Value = $vals[5]
    $ws.Range("S$row").Value = $vals[6]
}
